# Update power transfer values on the "aggr_exchange" sheet
# per commit: "Changed name of gen_trip and found a base for the 2030 scenario"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aggr_exchange")

# FI-EE: 963 -> 750
$ws.Range("B2").Value = 750

# NO_2-GB: -400 -> 0
$ws.Range("B6").Value = 0

# NO_2-NL: -700 -> -300
$ws.Range("B7").Value = -300

# SE_4-LT: 700 -> -500
$ws.Range("B11").Value = -500
